$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "-test1" suffixed usernames/emails to "-test"
# (names first, then emails, to mirror the original shared-string ordering)
for ($r = 2; $r -le 6; $r++) {
    $nameCell = $ws.Cells.Item($r, 1)
    $name = [string]$nameCell.Value2
    $nameCell.Value = $name.Replace("-test1", "-test")
}

for ($r = 2; $r -le 6; $r++) {
    $emailCell = $ws.Cells.Item($r, 3)
    $email = [string]$emailCell.Value2
    $emailCell.Value = $email.Replace("-test1", "-test")
}

# Add mailto hyperlinks on the email column (C2:C6), which also applies
# the built-in "Hyperlink" cell style to those cells.
for ($r = 2; $r -le 6; $r++) {
    $emailCell = $ws.Cells.Item($r, 3)
    $email = [string]$emailCell.Value2
    $ws.Hyperlinks.Add($emailCell, "mailto:" + $email)
}

# Update the active selection shown when the sheet was saved
$ws.Range("D13").Select()
